$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 7840
$ws.Range("I20").Value = 7840
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7840
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -7610
$ws.Range("N20").ClearContents()

$ws.Range("H26").Value = 7126.1333
$ws.Range("I26").Value = 6945.4546
$ws.Range("J26").Value = 7623
$ws.Range("K26").Value = 6945.4546
$ws.Range("L26").Value = 7623
$ws.Range("M26").Value = -6601.4546
$ws.Range("N26").Value = -8311

$ws.Range("H35").Value = 7840
$ws.Range("I35").Value = 7840
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 7840
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -7461
$ws.Range("N35").ClearContents()

$ws.Range("H48").Value = 12608.25
$ws.Range("I48").Value = 15017
$ws.Range("J48").Value = 12264.143
$ws.Range("K48").Value = 45051
$ws.Range("L48").Value = 36792.429
$ws.Range("M48").Value = -44759
$ws.Range("N48").Value = -37376.429

$ws.Range("H51").Value = 10395.094
$ws.Range("I51").Value = 8968
$ws.Range("J51").Value = 10870.792
$ws.Range("K51").Value = 8968
$ws.Range("L51").Value = 10870.792
$ws.Range("M51").Value = -8484
$ws.Range("N51").Value = -11838.792

$ws.Range("H56").Value = 12608.25
$ws.Range("I56").Value = 15017
$ws.Range("J56").Value = 12264.143
$ws.Range("K56").Value = 45051
$ws.Range("L56").Value = 36792.429
$ws.Range("M56").Value = -44517
$ws.Range("N56").Value = -37860.429

$ws.Range("H62").Value = 33338260
$ws.Range("I62").Value = 33338260
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 33338260
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -33337636

$ws.Range("H65").Value = 33338260
$ws.Range("I65").Value = 33338260
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 166691300
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -166688180

$ws.Range("H92").Value = 870.7
$ws.Range("I92").Value = 1020.58826
$ws.Range("J92").Value = 21.333334
$ws.Range("K92").Value = 1020.58826
$ws.Range("L92").Value = 21.333334
$ws.Range("M92").Value = 227.41174
$ws.Range("N92").Value = -2517.333334

$ws.Range("H112").Value = 2697.257
$ws.Range("I112").Value = 1985
$ws.Range("J112").Value = 2740.4243
$ws.Range("K112").Value = 5955
$ws.Range("L112").Value = 8221.2729
$ws.Range("M112").Value = -4847
$ws.Range("N112").Value = -10437.2729

$ws.Range("H116").Value = 4110.5
$ws.Range("I116").Value = 4089.4443
$ws.Range("J116").Value = 4300
$ws.Range("K116").Value = 4089.4443
$ws.Range("L116").Value = 4300
$ws.Range("M116").Value = -647.4443000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 36489.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 36489.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 36489.5
$ws.Range("N49").Value = -37009.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H46").Value = 20273.53
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 20273.53
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 20273.53
$ws.Range("N46").Value = -20869.53

$ws.Range("H105").Value = 38472868
$ws.Range("I105").Value = 41678732
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 41678732
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -41676985

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3803.3076
$ws.Range("I16").Value = 2976.25
$ws.Range("J16").Value = 5126.6
$ws.Range("K16").Value = 2976.25
$ws.Range("L16").Value = 5126.6
$ws.Range("M16").Value = -2689.25
$ws.Range("N16").Value = -5700.6

$ws.Range("H113").Value = 3803.3076
$ws.Range("I113").Value = 2976.25
$ws.Range("J113").Value = 5126.6
$ws.Range("K113").Value = 2976.25
$ws.Range("L113").Value = 5126.6
$ws.Range("M113").Value = -806.25
$ws.Range("N113").Value = -9466.6

$ws.Range("H132").Value = 47332.9
$ws.Range("I132").Value = 4462.7896
$ws.Range("J132").Value = 128786.1
$ws.Range("K132").Value = 13388.3688
$ws.Range("L132").Value = 386358.3
$ws.Range("M132").Value = -10858.3688
$ws.Range("N132").Value = -391418.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2831.9167
$ws.Range("I12").Value = 2992.3333
$ws.Range("J12").Value = 2778.4443
$ws.Range("K12").Value = 8976.999899999999
$ws.Range("L12").Value = 8335.332900000001
$ws.Range("M12").Value = -8803.999899999999
$ws.Range("N12").Value = -8681.332900000001

$ws.Range("H21").Value = 378.42856
$ws.Range("I21").Value = 100
$ws.Range("J21").Value = 489.8
$ws.Range("K21").Value = 300
$ws.Range("L21").Value = 1469.4
$ws.Range("M21").Value = -127
$ws.Range("N21").Value = -1815.4

$ws.Range("H35").Value = 3387
$ws.Range("I35").Value = 436.33334
$ws.Range("J35").Value = 5600
$ws.Range("K35").Value = 1309.00002
$ws.Range("L35").Value = 16800
$ws.Range("M35").Value = -1021.00002
$ws.Range("N35").Value = -17376

$ws.Range("H107").Value = 1180.1489
$ws.Range("I107").Value = 320.42856
$ws.Range("J107").Value = 1330.6
$ws.Range("K107").Value = 961.28568
$ws.Range("L107").Value = 3991.8
$ws.Range("M107").Value = 958.71432
$ws.Range("N107").Value = -7831.799999999999

$ws.Range("H119").Value = 399
$ws.Range("I119").Value = 399
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 1197
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 3641

$ws.Range("H122").Value = 7032.385
$ws.Range("I122").Value = 14683
$ws.Range("J122").Value = 2250.75
$ws.Range("K122").Value = 132147
$ws.Range("L122").Value = 20256.75
$ws.Range("M122").Value = -129697
$ws.Range("N122").Value = -25156.75

$ws.Range("H139").Value = 1296.625
$ws.Range("I139").Value = 1296.625
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3889.875
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1250.125

$ws.Range("H141").Value = 4276.8
$ws.Range("I141").Value = 4276.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 12830.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7650.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 48207.383
$ws.Range("I2").Value = 933.0833
$ws.Range("J2").Value = 111239.78
$ws.Range("K2").Value = 933.0833
$ws.Range("L2").Value = 111239.78
$ws.Range("M2").Value = -820.0833

$ws.Range("H35").Value = 13048.333
$ws.Range("I35").Value = 16072.5
$ws.Range("J35").Value = 7000
$ws.Range("K35").Value = 16072.5
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = -15774.5
$ws.Range("N35").Value = -7596

$ws.Range("H36").Value = 3166.6667
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = -1015
$ws.Range("N36").Value = -4970

$ws.Range("H46").Value = 33126.89
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 39734.57
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 39734.57
$ws.Range("M46").Value = -9844
$ws.Range("N46").Value = -40046.57

$ws.Range("H70").Value = 14916.667
$ws.Range("I70").Value = 13500.857
$ws.Range("J70").Value = 17748.285
$ws.Range("K70").Value = 13500.857
$ws.Range("L70").Value = 17748.285
$ws.Range("M70").Value = -13230.857

$ws.Range("H73").Value = 14916.667
$ws.Range("I73").Value = 13500.857
$ws.Range("J73").Value = 17748.285
$ws.Range("K73").Value = 13500.857
$ws.Range("L73").Value = 17748.285
$ws.Range("M73").Value = -12564.857

$ws.Range("H122").Value = 6562.1763
$ws.Range("I122").Value = 7405.7
$ws.Range("J122").Value = 5357.143
$ws.Range("K122").Value = 22217.1
$ws.Range("L122").Value = 16071.429
$ws.Range("M122").Value = -19767.1

$ws.Range("H132").Value = 7789.143
$ws.Range("I132").Value = 6354.143
$ws.Range("J132").Value = 9224.143
$ws.Range("K132").Value = 19062.429
$ws.Range("L132").Value = 27672.429
$ws.Range("M132").Value = -16532.429
$ws.Range("N132").Value = -32732.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5590

$ws.Range("H31").Value = 4717.6665
$ws.Range("I31").Value = 3879
$ws.Range("J31").Value = 5137
$ws.Range("K31").Value = 3879
$ws.Range("L31").Value = 5137
$ws.Range("M31").Value = -3631
$ws.Range("N31").Value = -5633

$ws.Range("H35").Value = 983.5
$ws.Range("I35").Value = 999.5
$ws.Range("J35").Value = 967.5
$ws.Range("K35").Value = 999.5
$ws.Range("L35").Value = 967.5
$ws.Range("M35").Value = -663.5
$ws.Range("N35").Value = -1639.5

$ws.Range("H61").Value = 8282.714
$ws.Range("I61").Value = 8282.714
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8282.714
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -8080.714

$ws.Range("H109").Value = 59285
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59285
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59285
$ws.Range("N109").Value = -62059

$ws.Range("H113").Value = 8282.714
$ws.Range("I113").Value = 8282.714
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8282.714
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6112.714

$ws.Range("H122").Value = 6219.3
$ws.Range("I122").Value = 4070.5715
$ws.Range("J122").Value = 11233
$ws.Range("K122").Value = 12211.7145
$ws.Range("L122").Value = 33699
$ws.Range("M122").Value = -9761.7145
$ws.Range("N122").Value = -38599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 250963.5
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 1284.6666
$ws.Range("K8").Value = 1000000
$ws.Range("L8").Value = 1284.6666
$ws.Range("M8").Value = -999860
$ws.Range("N8").Value = -1564.6666

$ws.Range("H28").Value = 39499.75
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 39499.75
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 39499.75
$ws.Range("N28").Value = -40195.75

$ws.Range("H33").Value = 29249.25
$ws.Range("I33").Value = 25000
$ws.Range("J33").Value = 30665.666
$ws.Range("K33").Value = 25000
$ws.Range("L33").Value = 30665.666
$ws.Range("M33").Value = -24750
$ws.Range("N33").Value = -31165.666

$ws.Range("H36").Value = 29249.25
$ws.Range("I36").Value = 25000
$ws.Range("J36").Value = 30665.666
$ws.Range("K36").Value = 25000
$ws.Range("L36").Value = 30665.666
$ws.Range("M36").Value = -24750
$ws.Range("N36").Value = -31165.666

$ws.Range("H38").Value = 16985.334
$ws.Range("I38").Value = 22028
$ws.Range("J38").Value = 6900
$ws.Range("K38").Value = 22028
$ws.Range("L38").Value = 6900
$ws.Range("M38").Value = -21555
$ws.Range("N38").Value = -7846

$ws.Range("H49").Value = 51243.832
$ws.Range("I49").Value = 33248.5
$ws.Range("J49").Value = 87234.5
$ws.Range("K49").Value = 33248.5
$ws.Range("L49").Value = 87234.5
$ws.Range("M49").Value = -33018.5
$ws.Range("N49").Value = -87694.5

$ws.Range("H52").Value = 15703.875
$ws.Range("I52").Value = 7673.6665
$ws.Range("J52").Value = 39794.5
$ws.Range("K52").Value = 7673.6665
$ws.Range("L52").Value = 39794.5
$ws.Range("M52").Value = -7447.6665
$ws.Range("N52").Value = -40246.5

$ws.Range("H132").Value = 6207.607
$ws.Range("I132").Value = 2378.923
$ws.Range("J132").Value = 9525.799999999999
$ws.Range("K132").Value = 7136.768999999999
$ws.Range("L132").Value = 28577.4
$ws.Range("M132").Value = -4606.768999999999
